# "Fruta / hortaliza, semanal" — add the week's new price-report row for
# Feria Lagunitas de Puerto Montt / Piña "Caramelo" / Segunda.
#
# The new observation is inserted as row 411 (sorted where it belongs in
# this sheet's date ordering), which pushes the existing rows 411-442 down
# to 412-443 — their contents are untouched, only their row numbers shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 411; everything at/after 411 shifts down by one.
$ws.Rows(411).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A411").Value = 4
$ws.Range("B411").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C411").Value = "Los Lagos"
$ws.Range("D411").Value = 45132
$ws.Range("E411").Value = 10
$ws.Range("F411").Value = "Fruta"
$ws.Range("G411").Value = 100108
$ws.Range("H411").Value = "Tropicales y subtropicales"
$ws.Range("I411").Value = 100108005
$ws.Range("J411").Value = "Piña"
$ws.Range("K411").Value = "Caramelo"
$ws.Range("L411").Value = "Segunda"
$ws.Range("M411").Value = 120
$ws.Range("N411").Value = 22000
$ws.Range("O411").Value = 22000
$ws.Range("P411").Value = 22000
$ws.Range("Q411").Value = "$/caja 14 unidades"
$ws.Range("R411").Value = "Ecuador"
$ws.Range("S411").Value = 1571
$ws.Range("T411").Value = 14
